$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Insert 5 fresh rows above the existing two data rows (rows 1-2 become rows 6-7).
$ws.Rows("1:5").Insert()

# Helper-style approach: for each row, force the target range to Text format so
# numeric-looking values ("196.63", "3305014528", ...) are stored as shared
# strings (t="s") instead of being coerced to numeric cells, then strip the
# number-format style back off so the cells end up on the default style (no
# "s" attribute), matching the rest of the sheet.
function Set-RowValues($rowNum, $values) {
    $lastCol = [char](65 + $values.Count - 1)
    $addr = "A" + $rowNum + ":" + $lastCol + $rowNum
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    for ($i = 0; $i -lt $values.Count; $i++) {
        $col = [string]([char](65 + $i))
        $cellAddr = $col + $rowNum
        $ws.Range($cellAddr).Value = $values[$i]
    }
    $rng.Style = "Normal"
}

Set-RowValues 1 @("iuliia.1", "CAD", "196.63", "14", "Visa", "3305014528")
Set-RowValues 2 @("iuliia.4", "AED", "254.6", "40", "MasterCard", "3398670848")
Set-RowValues 3 @("iuliia.4", "KWD", "24.26", "4", "American Express", "3321217280")
Set-RowValues 4 @("iuliia.2", "NZD", "837.4", "200", "Visa", "3316641632")
Set-RowValues 5 @("iuliia.6", "AUD", "1562.88", "400", "MasterCard", "3369803520")

# Rows 6-7 already carry the old row-1/row-2 data (shifted down by the insert);
# only their Amount (C) and Transaction-Id (F) values actually change.
Set-RowValues 6 @("iuliia.6", "GBP", "680.82", "300", "American Express", "3300547584")
Set-RowValues 7 @("iuliia.1", "CAD", "210.63", "14", "Visa")

Write-Output "done"
